$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("H62").Value = 8567.666999999999
$ws.Range("I62").Value = 2204.6667
$ws.Range("J62").Value = 11749.167
$ws.Range("K62").Value = 2204.6667
$ws.Range("L62").Value = 11749.167
$ws.Range("M62").Value = -1580.6667
$ws.Range("N62").Value = -12997.167
$ws.Range("H65").Value = 8567.666999999999
$ws.Range("I65").Value = 2204.6667
$ws.Range("J65").Value = 11749.167
$ws.Range("K65").Value = 11023.3335
$ws.Range("L65").Value = 58745.835
$ws.Range("M65").Value = -7903.333500000001
$ws.Range("N65").Value = -64985.835
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("H98").Value = 784.2353000000001
$ws.Range("I98").Value = 739.5454999999999
$ws.Range("J98").Value = 866.1667
$ws.Range("K98").Value = 739.5454999999999
$ws.Range("L98").Value = 866.1667
$ws.Range("M98").Value = 758.4545000000001
$ws.Range("N98").Value = -3862.1667
$ws.Range("H122").Value = 784.2353000000001
$ws.Range("I122").Value = 739.5454999999999
$ws.Range("J122").Value = 866.1667
$ws.Range("K122").Value = 2218.6365
$ws.Range("L122").Value = 2598.5001
$ws.Range("M122").Value = 231.3635000000004
$ws.Range("N122").Value = -7498.5001
$ws.Range("H125").Value = 3738.5
$ws.Range("I125").Value = 3738.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 33646.5
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = -31186.5
$ws.Range("M125").ClearContents()
$ws.Range("H129").Value = 1578.1538
$ws.Range("I129").Value = 968.5
$ws.Range("J129").Value = 2553.6
$ws.Range("K129").Value = 2905.5
$ws.Range("L129").Value = 7660.799999999999
$ws.Range("M129").Value = 2094.5
$ws.Range("N129").Value = -17660.8
$ws.Range("H131").Value = 2157.7693
$ws.Range("I131").Value = 2096.3333
$ws.Range("J131").Value = 2895
$ws.Range("K131").Value = 6288.999899999999
$ws.Range("L131").Value = 8685
$ws.Range("M131").Value = -1248.999899999999
$ws.Range("N131").Value = -18765

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4105.95
$ws.Range("I74").Value = 4075
$ws.Range("J74").Value = 4198.8
$ws.Range("K74").Value = 4075
$ws.Range("L74").Value = 4198.8
$ws.Range("M74").Value = -3201
$ws.Range("N74").Value = -5946.8
$ws.Range("H77").Value = 4105.95
$ws.Range("I77").Value = 4075
$ws.Range("J77").Value = 4198.8
$ws.Range("K77").Value = 20375
$ws.Range("L77").Value = 20994
$ws.Range("M77").Value = -16007
$ws.Range("N77").Value = -29730
$ws.Range("H132").Value = 1917.5294
$ws.Range("I132").Value = 1909.5483
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5728.644899999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3198.644899999999
$ws.Range("N132").Value = -11060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3650.6
$ws.Range("I134").Value = 1128.9565
$ws.Range("J134").Value = 32649.5
$ws.Range("K134").Value = 3386.8695
$ws.Range("L134").Value = 97948.5
$ws.Range("M134").Value = -851.8694999999998
$ws.Range("N134").Value = -103018.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.933334
$ws.Range("I7").Value = 51.153847
$ws.Range("J7").Value = 199.5
$ws.Range("K7").Value = 51.153847
$ws.Range("L7").Value = 199.5
$ws.Range("M7").Value = 61.846153
$ws.Range("N7").Value = -425.5
$ws.Range("H22").Value = 2244.6667
$ws.Range("I22").Value = 1601.75
$ws.Range("J22").Value = 2759
$ws.Range("K22").Value = 1601.75
$ws.Range("L22").Value = 2759
$ws.Range("M22").Value = -1251.75
$ws.Range("N22").Value = -3459
$ws.Range("H132").Value = 3075.8518
$ws.Range("I132").Value = 2830.348
$ws.Range("J132").Value = 4487.5
$ws.Range("K132").Value = 8491.044
$ws.Range("L132").Value = 13462.5
$ws.Range("M132").Value = -5961.044
$ws.Range("N132").Value = -18522.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2164.5264
$ws.Range("I131").Value = 1772.125
$ws.Range("J131").Value = 2449.9092
$ws.Range("K131").Value = 5316.375
$ws.Range("L131").Value = 7349.7276
$ws.Range("M131").Value = -276.375
$ws.Range("N131").Value = -17429.7276
$ws.Range("H138").Value = 4025.9
$ws.Range("I138").Value = 1193.1666
$ws.Range("J138").Value = 8275
$ws.Range("K138").Value = 3579.4998
$ws.Range("L138").Value = 24825
$ws.Range("M138").Value = 1560.5002
$ws.Range("N138").Value = -35105

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 111151.22
$ws.Range("I2").Value = 250033.75
$ws.Range("J2").Value = 45.2
$ws.Range("K2").Value = 250033.75
$ws.Range("L2").Value = 45.2
$ws.Range("M2").Value = -249920.75
$ws.Range("N2").Value = -271.2
$ws.Range("H62").Value = 2500000
$ws.Range("I62").Value = 2500000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2500000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2499314
$ws.Range("H65").Value = 2500000
$ws.Range("I65").Value = 2500000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7500000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -7496568
$ws.Range("H122").Value = 1469.0769
$ws.Range("I122").Value = 1469.0769
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4407.2307
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1957.2307

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2790.125
$ws.Range("I7").Value = 2337.25
$ws.Range("J7").Value = 4148.75
$ws.Range("K7").Value = 2337.25
$ws.Range("L7").Value = 4148.75
$ws.Range("M7").Value = -2225.25
$ws.Range("N7").Value = -4372.75
$ws.Range("H16").Value = 3324
$ws.Range("I16").Value = 3453.7
$ws.Range("J16").Value = 2999.75
$ws.Range("K16").Value = 3453.7
$ws.Range("L16").Value = 2999.75
$ws.Range("M16").Value = -3283.7
$ws.Range("N16").Value = -3339.75
$ws.Range("H22").Value = 833
$ws.Range("I22").Value = 749.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 749.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -454.5
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 833
$ws.Range("I27").Value = 749.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 749.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -642.5
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 4672
$ws.Range("I40").Value = 3783.1667
$ws.Range("J40").Value = 10005
$ws.Range("K40").Value = 3783.1667
$ws.Range("L40").Value = 10005
$ws.Range("M40").Value = -3647.1667
$ws.Range("N40").Value = -10277
$ws.Range("H46").Value = 2858.65
$ws.Range("I46").Value = 1322.125
$ws.Range("J46").Value = 3883
$ws.Range("K46").Value = 1322.125
$ws.Range("L46").Value = 3883
$ws.Range("M46").Value = -1134.125
$ws.Range("N46").Value = -4259
$ws.Range("H55").Value = 812.05
$ws.Range("I55").Value = 928.2857
$ws.Range("J55").Value = 540.8333
$ws.Range("K55").Value = 928.2857
$ws.Range("L55").Value = 540.8333
$ws.Range("M55").Value = -755.2857
$ws.Range("N55").Value = -886.8333
$ws.Range("H68").Value = 8500.333000000001
$ws.Range("I68").Value = 6001
$ws.Range("J68").Value = 9750
$ws.Range("K68").Value = 6001
$ws.Range("L68").Value = 9750
$ws.Range("M68").Value = -5252
$ws.Range("N68").Value = -11248
$ws.Range("H71").Value = 8500.333000000001
$ws.Range("I71").Value = 6001
$ws.Range("J71").Value = 9750
$ws.Range("K71").Value = 30005
$ws.Range("L71").Value = 48750
$ws.Range("M71").Value = -26261
$ws.Range("N71").Value = -56238
$ws.Range("H122").Value = 3263.2144
$ws.Range("I122").Value = 3415.111
$ws.Range("J122").Value = 2989.8
$ws.Range("K122").Value = 10245.333
$ws.Range("L122").Value = 8969.400000000001
$ws.Range("M122").Value = -7795.332999999999
$ws.Range("N122").Value = -13869.4
$ws.Range("H126").Value = 2790.125
$ws.Range("I126").Value = 2337.25
$ws.Range("J126").Value = 4148.75
$ws.Range("K126").Value = 7011.75
$ws.Range("L126").Value = 12446.25
$ws.Range("M126").Value = -4541.75
$ws.Range("N126").Value = -17386.25
$ws.Range("H132").Value = 5174.2856
$ws.Range("I132").Value = 3951.0715
$ws.Range("J132").Value = 7620.7144
$ws.Range("K132").Value = 11853.2145
$ws.Range("L132").Value = 22862.1432
$ws.Range("M132").Value = -9323.2145
$ws.Range("N132").Value = -27922.1432
$ws.Range("H136").Value = 3268.4546
$ws.Range("I136").Value = 1900.75
$ws.Range("J136").Value = 4050
$ws.Range("K136").Value = 5702.25
$ws.Range("L136").Value = 12150
$ws.Range("M136").Value = -3152.25
$ws.Range("N136").Value = -17250

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5719.8
$ws.Range("I122").Value = 4639.8
$ws.Range("J122").Value = 6799.8
$ws.Range("K122").Value = 13919.4
$ws.Range("L122").Value = 20399.4
$ws.Range("M122").Value = -11469.4
$ws.Range("N122").Value = -25299.4
